$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 391.25
$ws.Range("I33").Value = 393.64285
$ws.Range("J33").Value = 374.5
$ws.Range("K33").Value = 393.64285
$ws.Range("L33").Value = 374.5
$ws.Range("M33").Value = -164.64285
$ws.Range("N33").Value = -832.5
# Row 106
$ws.Range("H106").Value = 11564.462
$ws.Range("I106").Value = 11861.833
$ws.Range("K106").Value = 11861.833
$ws.Range("M106").Value = -11230.833
# Row 129
$ws.Range("H129").Value = 2333.3125
$ws.Range("I129").Value = 1396.375
$ws.Range("J129").Value = 3270.25
$ws.Range("K129").Value = 4189.125
$ws.Range("L129").Value = 9810.75
$ws.Range("M129").Value = 810.875
$ws.Range("N129").Value = -19810.75
# Row 132
$ws.Range("H132").Value = 4419.893
$ws.Range("I132").Value = 4579.143
$ws.Range("K132").Value = 13737.429
$ws.Range("M132").Value = -11207.429
# Row 137
$ws.Range("H137").Value = 874.1429000000001
$ws.Range("I137").Value = 769.8333
$ws.Range("K137").Value = 2309.4999
$ws.Range("M137").Value = 240.5001000000002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 342.8
$ws.Range("I5").Value = 379
$ws.Range("K5").Value = 379
$ws.Range("M5").Value = -267
# Row 32
$ws.Range("H32").Value = 5071.4414
$ws.Range("I32").Value = 5071.4414
$ws.Range("K32").Value = 5071.4414
$ws.Range("M32").Value = -4784.4414
# Row 61
$ws.Range("H61").Value = 3469.08
$ws.Range("I61").Value = 2621.4
$ws.Range("J61").Value = 4740.6
$ws.Range("K61").Value = 2621.4
$ws.Range("L61").Value = 4740.6
$ws.Range("M61").Value = -2409.4
$ws.Range("N61").Value = -5164.6
# Row 97
$ws.Range("H97").Value = 2124
$ws.Range("I97").Value = 1073.9333
$ws.Range("K97").Value = 1073.9333
$ws.Range("M97").Value = -577.9332999999999
# Row 122
$ws.Range("H122").Value = 2026
$ws.Range("I122").Value = 1839.8
$ws.Range("K122").Value = 5519.4
$ws.Range("M122").Value = -3069.4
# Row 136
$ws.Range("H136").Value = 3469.08
$ws.Range("I136").Value = 2621.4
$ws.Range("J136").Value = 4740.6
$ws.Range("K136").Value = 7864.200000000001
$ws.Range("L136").Value = 14221.8
$ws.Range("M136").Value = -5314.200000000001
$ws.Range("N136").Value = -19321.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 342.8
$ws.Range("I4").Value = 379
$ws.Range("K4").Value = 379
$ws.Range("M4").Value = -264
# Row 94
$ws.Range("H94").Value = 1000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1000
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -1902

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Range("H12").Value = 2929.8
$ws.Range("J12").Value = 2999.8333
$ws.Range("L12").Value = 2999.8333
$ws.Range("N12").Value = -3339.8333
# Row 16
$ws.Range("H16").Value = 6734.75
$ws.Range("I16").Value = 470
$ws.Range("J16").Value = 12999.5
$ws.Range("K16").Value = 470
$ws.Range("L16").Value = 12999.5
$ws.Range("M16").Value = -183
$ws.Range("N16").Value = -13573.5
# Row 22
$ws.Range("H22").Value = 3078483.2
$ws.Range("I22").Value = 1555.25
$ws.Range("J22").Value = 8001568
$ws.Range("K22").Value = 1555.25
$ws.Range("L22").Value = 8001568
$ws.Range("M22").Value = -1205.25
$ws.Range("N22").Value = -8002268
# Row 69
$ws.Range("H69").Value = 14029.333
$ws.Range("I69").Value = 14029.333
$ws.Range("K69").Value = 14029.333
$ws.Range("M69").Value = -13280.333
# Row 72
$ws.Range("H72").Value = 14029.333
$ws.Range("I72").Value = 14029.333
$ws.Range("K72").Value = 42087.999
$ws.Range("M72").Value = -38343.999
# Row 113
$ws.Range("H113").Value = 6734.75
$ws.Range("I113").Value = 470
$ws.Range("J113").Value = 12999.5
$ws.Range("K113").Value = 470
$ws.Range("L113").Value = 12999.5
$ws.Range("M113").Value = 1700
$ws.Range("N113").Value = -17339.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 8035.0527
$ws.Range("I56").Value = 8035.0527
$ws.Range("K56").Value = 8035.0527
$ws.Range("M56").Value = -7505.0527
# Row 68
$ws.Range("H68").Value = 1200
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1200
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 3600
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -5222
# Row 71
$ws.Range("H71").Value = 1200
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1200
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 10800
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -18912
# Row 119
$ws.Range("H119").Value = 1999.3334
$ws.Range("I119").Value = 1999.3334
$ws.Range("K119").Value = 5998.0002
$ws.Range("M119").Value = -1160.0002
# Row 122
$ws.Range("H122").Value = 1948
$ws.Range("I122").Value = 1898
$ws.Range("K122").Value = 17082
$ws.Range("M122").Value = -14632

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Range("H24").Value = 111518.82
$ws.Range("J24").Value = 22670.7
$ws.Range("L24").Value = 22670.7
$ws.Range("N24").Value = -23016.7
# Row 126
$ws.Range("H126").Value = 4466.3335
$ws.Range("I126").Value = 4199.5
$ws.Range("K126").Value = 12598.5
$ws.Range("M126").Value = -10128.5
# Row 132
$ws.Range("H132").Value = 2267.65
$ws.Range("I132").Value = 2236.2144
$ws.Range("K132").Value = 6708.6432
$ws.Range("M132").Value = -4178.6432

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 17
$ws.Range("H17").Value = 9249.5
$ws.Range("J17").Value = 9249.5
$ws.Range("L17").Value = 9249.5
$ws.Range("N17").Value = -9589.5
# Row 25
$ws.Range("H25").Value = 17671.666
$ws.Range("I25").Value = 8007
$ws.Range("K25").Value = 8007
$ws.Range("M25").Value = -7777
# Row 46
$ws.Range("H46").Value = 2388.2222
$ws.Range("J46").Value = 3415
$ws.Range("L46").Value = 3415
$ws.Range("N46").Value = -3791
# Row 61
$ws.Range("H61").Value = 3000
$ws.Range("I61").Value = 3000
$ws.Range("K61").Value = 3000
$ws.Range("M61").Value = -2798
# Row 113
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
# Row 122
$ws.Range("H122").Value = 4051.75
$ws.Range("I122").Value = 3981.8
$ws.Range("J122").Value = 4168.3335
$ws.Range("K122").Value = 11945.4
$ws.Range("L122").Value = 12505.0005
$ws.Range("M122").Value = -9495.400000000001
$ws.Range("N122").Value = -17405.0005
# Row 132
$ws.Range("H132").Value = 2259.625
$ws.Range("I132").Value = 1661.6
$ws.Range("K132").Value = 4984.799999999999
$ws.Range("M132").Value = -2454.799999999999
